# Insert a new weekly price-report row at row 198 (pushing the existing
# rows 198:254 down to 199:255), then populate the new row with the new
# observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 198..254 down by one to make room for the new record.
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new data point.
$ws.Cells.Item(198, 1).Value = 1
$ws.Cells.Item(198, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(198, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(198, 4).Value = 44642
$ws.Cells.Item(198, 5).Value = 15
$ws.Cells.Item(198, 6).Value = 100114013
$ws.Cells.Item(198, 7).Value = "Zanahoria"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 70
$ws.Cells.Item(198, 11).Value = 24000
$ws.Cells.Item(198, 12).Value = 25000
$ws.Cells.Item(198, 13).Value = 24500
$ws.Cells.Item(198, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(198, 15).Value = "Valle de Camiña"
$ws.Cells.Item(198, 16).Value = 980
$ws.Cells.Item(198, 17).Value = 25
$ws.Cells.Item(198, 18).Value = "Hortaliza"
